# Insert a new data row above the current row 183 (shifts the existing
# rows 183-235 down to 184-236, matching the target diff where:
#   - dimension grows from A1:T235 to A1:T236
#   - every row from 184..236 now carries the values previously held by
#     the row one above it (183..235)
#   - the brand-new row 183 is a copy of the old row 183 except for the
#     Fecha (D) and Volumen (M) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 183..235 down to 184..236, leaving a blank row 183 behind.
$ws.Rows.Item(183).Insert()

# Populate the new row 183 with the same record as the (now shifted)
# row 184, except for the date and volume columns.
$ws.Cells.Item(183, 1).Value  = 4
$ws.Cells.Item(183, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(183, 3).Value  = "Los Lagos"
$ws.Cells.Item(183, 4).Value  = 44663
$ws.Cells.Item(183, 5).Value  = 10
$ws.Cells.Item(183, 6).Value  = "Fruta"
$ws.Cells.Item(183, 7).Value  = 100101
$ws.Cells.Item(183, 8).Value  = "Berries"
$ws.Cells.Item(183, 9).Value  = 100112025
$ws.Cells.Item(183, 10).Value = "Frutilla"
$ws.Cells.Item(183, 11).Value = "Sin especificar"
$ws.Cells.Item(183, 12).Value = "Primera"
$ws.Cells.Item(183, 13).Value = 500
$ws.Cells.Item(183, 14).Value = 8500
$ws.Cells.Item(183, 15).Value = 9000
$ws.Cells.Item(183, 16).Value = 8750
$ws.Cells.Item(183, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(183, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(183, 19).Value = 1250
$ws.Cells.Item(183, 20).Value = 7
